# Add 5 new algorithm problems (about dp / array / queue / string) to the
# tracking sheet, continuing directly after the existing last row (161).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- new row data: [A number, B title, C topic, D difficulty, hyperlink-or-$null]
$newRows = @(
    @{ Row = 162; Num = 163; Title = "剑指 Offer 63. 股票的最大利润";     Topic = "动态规划"; Level = "中";  Link = "https://leetcode.cn/problems/gu-piao-de-zui-da-li-run-lcof/" },
    @{ Row = 163; Num = 164; Title = "面试题13. 机器人的运动范围";       Topic = "队列";     Level = "中";  Link = "https://leetcode.cn/problems/ji-qi-ren-de-yun-dong-fan-wei-lcof/" },
    @{ Row = 164; Num = 165; Title = "面试题45. 把数组排成最小的数";     Topic = "字符串";   Level = "中";  Link = "https://leetcode.cn/problems/ba-shu-zu-pai-cheng-zui-xiao-de-shu-lcof/" },
    @{ Row = 165; Num = 166; Title = "面试题59 - II. 队列的最大值";      Topic = "队列";     Level = "中";  Link = "https://leetcode.cn/problems/dui-lie-de-zui-da-zhi-lcof/" },
    @{ Row = 166; Num = 167; Title = "面试题61. 扑克牌中的顺子";         Topic = "数组";     Level = "简单"; Link = $null }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Num
    $ws.Cells.Item($r.Row, 2).Value = $r.Title
    $ws.Cells.Item($r.Row, 3).Value = $r.Topic
    $ws.Cells.Item($r.Row, 4).Value = $r.Level

    if ($r.Link) {
        $ws.Hyperlinks.Add($ws.Cells.Item($r.Row, 2), $r.Link) | Out-Null
    }

    # Match the look of the other "title" cells in column B (hyperlink-style
    # font, left aligned) regardless of whether a live hyperlink was added.
    $ws.Cells.Item($r.Row, 2).Style = "Hyperlink"
    $ws.Cells.Item($r.Row, 2).HorizontalAlignment = -4131
}

# Keep the active selection near where the author left off editing.
$ws.Range("C169").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 124
$excel.ActiveWindow.ScrollColumn = 1
